# RMA Complete Flow(Repair)-SO To RMA Receipt TO Replacement SO.xlsx
# "updated testcases- API-WOrk order"
#
# The "RMA Details Maintenance Grid" sheet holds three rows of RMA test
# fixture data (RMA number, Sales Order Line id, Shipper Line id). This
# edit refreshes that fixture to a newly-generated RMA ("ZH9H") in place
# of the previous one ("VXRY").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 - RMA-ZH9H-001
$ws.Range("E2").Value = "RMA-ZH9H-001"
$ws.Range("F2").Value = "RMA-ZH9H-1-1"
$ws.Range("J2").Value = "a7s5f000000xL38AAE"

# Row 3 - RMA-ZH9H-002
$ws.Range("E3").Value = "RMA-ZH9H-002"
$ws.Range("F3").Value = "RMA-ZH9H-1-2"
$ws.Range("J3").Value = "a7s5f000000xL39AAE"

# Row 4 - RMA-ZH9H-003
$ws.Range("E4").Value = "RMA-ZH9H-003"
$ws.Range("F4").Value = "RMA-ZH9H-1-3"
$ws.Range("J4").Value = "a7s5f000000xL3AAAU"
